# Append a new data row (row 7) to the active worksheet, mirroring the
# structure of the existing rows (A = date/time, B:M = numeric values,
# N = "Named" text label).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 7

$ws.Cells.Item($row, 1).Value = 42611.884664351855
$ws.Cells.Item($row, 2).Value = -34
$ws.Cells.Item($row, 3).Value = 42
$ws.Cells.Item($row, 4).Value = 57
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 99
$ws.Cells.Item($row, 7).Value = 13715
$ws.Cells.Item($row, 8).Value = 10482
$ws.Cells.Item($row, 9).Value = 608
$ws.Cells.Item($row, 10).Value = 64
$ws.Cells.Item($row, 11).Value = 87
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 6
$ws.Cells.Item($row, 14).Value = "Named"

# Column A keeps the same date-number-format style as the rows above it.
# Copy just the formatting (not the value) from the cell above so that the
# existing style entry is reused instead of a new style being created.
$ws.Cells.Item($row - 1, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
